$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 69 (KgType) - two cells flip from "No" to "Yes"
# ---------------------------------------------------------------------------
$ws.Range("Q69").Value = "Yes"
$ws.Range("AA69").Value = "Yes"

# ---------------------------------------------------------------------------
# Row 71 - new top-level type "KgTypeReference" (same pattern as row 69)
# ---------------------------------------------------------------------------
$cols = @("P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# New label strings are introduced to the shared-strings table in the same
# order the author must have typed them (label column first, in row
# 71, 77, 75, 73, 72, 74, 76 order) before the rest of the grid is filled in.
$ws.Range("C71").Value = "KgTypeReference"
$ws.Range("D77").Value = "KgUserType"
$ws.Range("D75").Value = "KgNullableType"
$ws.Range("D73").Value = "KgFunctionType"
$ws.Range("D72").Value = "KgDynamicType"
$ws.Range("D74").Value = "KgInferredType"
$ws.Range("D76").Value = "KgParenthesizedType"

foreach ($col in $cols) {
    $ws.Range($col + "71").Value = "No"
}
$ws.Range("AA71").Value = "Yes"

# ---------------------------------------------------------------------------
# Rows 72-77 - concrete subtypes of KgTypeReference, each marked with the
# literal quote character across the P:AC columns (same convention used by
# the existing KgUserType-style rows such as 66/67).
# ---------------------------------------------------------------------------
$subtypeRows = @(72, 73, 74, 75, 76, 77)

foreach ($r in $subtypeRows) {
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = '"'
    }
}

# ---------------------------------------------------------------------------
# Selection / view bookkeeping to mirror the author's saved state
# ---------------------------------------------------------------------------
$ws.Range("H80").Select()
